# Add three new "A_SEKTOREN_*" rows (A_SEKTOREN_25, A_SEKTOREN_51,
# A_SEKTOREN_55bc) right before the existing A_SEKTOREN_HOCHSCH row
# (originally row 73), pushing all subsequent rows down by three.
# Also fix a capitalisation typo that rides along in the same upload
# ("Technical compliance" -> "Technical Compliance" in column C of the
# row that carries A_SERIES_FATFTEC).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above row 73 (one at a time so each new row
# is truly blank rather than a single 3-row block insert).
$ws.Rows.Item(73).Insert()
$ws.Rows.Item(74).Insert()
$ws.Rows.Item(75).Insert()

# The row that used to be 73 (A_SEKTOREN_HOCHSCH) is now row 76 and
# still carries the correct cell formatting/style for this table; copy
# it onto the freshly inserted blank rows so the new rows look like
# the rest of the data table instead of Excel's default insert style.
$ws.Range("A76:D76").Copy()
$ws.Range("A73:D75").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the values for the three new rows.
$ws.Range("A73").Value = "A_SEKTOREN_25"
$ws.Range("B73").Value = "K_SEKTOREN"
$ws.Range("C73").Value = "Nachhaltige bewirtschaftete Fläche"
$ws.Range("D73").Value = "XXXNachhaltige bewirtschaftete Fläche"

$ws.Range("A74").Value = "A_SEKTOREN_51"
$ws.Range("B74").Value = "K_SEKTOREN"
$ws.Range("C74").Value = "Terrestrische Schutzgebiete und Mangrovenwälder in marinen Schutzgebieten"
$ws.Range("D74").Value = "XXXTerrestrische Schutzgebiete und Mangrovenwälder in marinen Schutzgebieten"

$ws.Range("A75").Value = "A_SEKTOREN_55bc"
$ws.Range("B75").Value = "K_SEKTOREN"
$ws.Range("C75").Value = "Waldfläche nachhaltige Bewirtschaftung und Wiederherstellung"
$ws.Range("D75").Value = "XXXWaldfläche nachhaltige Bewirtschaftung und Wiederherstellung"

# Small accompanying text fix: the row for A_SERIES_FATFTEC (now row
# 106 after the shift) gets its German label's capitalisation
# corrected from "compliance" to "Compliance".
$ws.Range("C106").Value = "<u>Infolinie:</u> Technical Compliance"
